$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data row for "南京耀多信息技术有限公司" (old row 3) was removed from the
# source list backing this sheet. Column A holds a static sequence number
# that is independent of row content, so it is left untouched; only the data
# columns B:P are pulled up from the row below for every row from the
# deleted entry's position down through the last data row. A leading
# apostrophe forces each value to land back in the cell as plain text (not a
# re-interpreted number/percentage), matching how the source data is stored.
# The now-duplicate last row is then removed outright so the sheet shrinks
# from 12 to 11 data rows (dimension A1:P13 -> A1:P12).
for ($r = 3; $r -le 12; $r++) {
    for ($c = 2; $c -le 16; $c++) {
        $src = $ws.Cells.Item($r + 1, $c)
        $dst = $ws.Cells.Item($r, $c)
        $dst.Value2 = "'" + $src.Value2
    }
}
$ws.Rows.Item(13).Delete()
